# Apply updated cryptocurrency price/volume data scraped on 2023-02-21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.711.02'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '1.680.10'
$ws.Range('E3').Value = '  -1.66%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.51'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3928'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3960'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.37%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.003'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.34%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '51.82'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.23%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.417'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -5.20%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08660'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '25.31'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.96%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.355'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.02%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.814'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.99%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001322'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.87%  '
$ws.Range('D17').Value = '1.640.34'
$ws.Range('E17').Value = '  -2.54%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '93.81'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.66%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.07104'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.41%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '20.23'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -4.95%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.132'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.50%  '
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.04'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.47%  '
$ws.Range('D24').Value = '24.717.66'
$ws.Range('E24').Value = '  -0.47%  '
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '23.67'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.52%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.783'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -7.54%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '162.53'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.44%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.811'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.26%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '151.23'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +4.06%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.853'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -7.40%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.381'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +5.49%  '
$ws.Range('D33').Value = '1.837.14'
$ws.Range('E33').Value = '  -4.04%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08475'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.07%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.03095'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.67%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.013'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.47%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.977'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.54%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2799'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.81%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.09509'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.97%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '10.54'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.84%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.7968'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.79%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.482'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '13.67'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.45%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.78'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -4.13%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.581'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -4.43%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.7154'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.58%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.173'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.78%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.08666'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +4.42%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.002'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.334'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.80%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '138.27'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.16%  '
